$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (UNH) - updated metrics
$ws.Range("D2").Value = 331.13
$ws.Range("E2").Value = 56.9
$ws.Range("F2").Value = 0.41
$ws.Range("K2").Value = 57.2
$ws.Range("N2").Value = 50.60178744571824

# Row 3 - now MetLife, Inc. (MET) instead of AIG, with updated metrics
$ws.Range("B3").Value = "MetLife, Inc."
$ws.Range("C3").Value = "MET"
$ws.Range("D3").Value = 78.42
$ws.Range("E3").Value = 48.9
$ws.Range("F3").Value = 2.43
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 54.4
$ws.Range("N3").Value = 50.60178744571824

# Row 4 - now American International Group, I (AIG) instead of MetLife, with updated metrics
$ws.Range("B4").Value = "American International Group, I"
$ws.Range("C4").Value = "AIG"
$ws.Range("D4").Value = 77.22
$ws.Range("E4").Value = 46.1
$ws.Range("F4").Value = 1.39
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 46
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 51.6
$ws.Range("N4").Value = 50.60178744571824

# Row 5 (PRU) - updated metrics
$ws.Range("D5").Value = 111.55
$ws.Range("E5").Value = 69.90000000000001
$ws.Range("F5").Value = 3.05
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 33
$ws.Range("K5").Value = 46.4
$ws.Range("N5").Value = 50.60178744571824
